$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing diagonal values (QoQ error -> re-based / recomputed values)
$ws.Range("K10").Value = 0.5403464745801891

$ws.Range("J11").Value = 0.4969491838668565
$ws.Range("K11").Value = 0.2970525035592049

$ws.Range("I12").Value = 0.5011245128056051
$ws.Range("J12").Value = 0.2858677898194339

$ws.Range("H13").Value = 0.4852787037784192
$ws.Range("I13").Value = 0.2775335613519331

$ws.Range("G14").Value = 0.4539510573947921
$ws.Range("H14").Value = 0.2743085116504074

$ws.Range("F15").Value = 0.4663391832225094
$ws.Range("G15").Value = 0.2534447081011285

$ws.Range("E16").Value = 0.4814444548743619
$ws.Range("F16").Value = 0.2766837437271186

$ws.Range("D17").Value = 0.4184715358843989
$ws.Range("E17").Value = 0.2867219094086165

$ws.Range("C18").Value = 0.5177895860664353
$ws.Range("D18").Value = 0.1751453671933744

$ws.Range("B19").Value = 0.5618492773058843
$ws.Range("C19").Value = 0.1965658720679752

$ws.Range("B20").Value = 0.4328090033804217
